# Daily attendance processing - 2025-12-17 03:16:09
# Reorders the "Recorded By" (column G) author list so that entries whose
# list ends with "System"/"system" have their comma-separated parts
# reversed (System moves to the front of the list).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($null -eq $val) { continue }
    if ($val -isnot [string]) { continue }

    $trimmed = $val.Trim()
    if ($trimmed.ToLower().EndsWith(", system")) {
        $parts = $val -split ","
        for ($i = 0; $i -lt $parts.Length; $i++) {
            $parts[$i] = $parts[$i].Trim()
        }
        $reversed = $parts[($parts.Length - 1)..0]
        $cell.Value = [string]::Join(", ", $reversed)
    }
}
